$d = $word.ActiveDocument

# Locate the paragraph that ends with "My concentration is databases", so the
# new paragraph can be inserted right after it.
$targetIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -match "databases") {
        $targetIdx = $i
    }
}

$targetPara = $d.Paragraphs($targetIdx)

# Add a new paragraph after it containing the professor's name.
$targetPara.Range.InsertParagraphAfter()
$newIdx = $targetIdx + 1
$newPara = $d.Paragraphs($newIdx)
$newPara.Range.InsertAfter("My professor is Peter WOLCOTT")

# The "_GoBack" bookmark used to sit at the end of the "...14, 2018" paragraph;
# remove it there since the most recent edit now happens in the new paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create "_GoBack" as a zero-length bookmark right after the new text, i.e.
# at the end of the new paragraph's content (just before its paragraph mark).
$newPara = $d.Paragraphs($newIdx)
$insertPos = $newPara.Range.End - 1

# Bookmarking an empty (collapsed) range directly at a paragraph's final
# character position anchors it incorrectly in this runtime, so briefly
# insert a marker character to anchor a real (non-empty) range, bookmark
# that range, then delete the marker - the bookmark collapses back down
# and stays anchored at the same spot.
$marker = $d.Range($insertPos, $insertPos)
$marker.InsertAfter("X")
$markerRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$d.Range($insertPos, $insertPos + 1).Delete()
